# Korrektur Testcode und Anpassung UML Diagramm
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Simple Test Suite")

# --- Row 9: Eingabewert -3 statt Text "3.5" ---
$ws.Range("E9").Value = -3

# --- Row 10: Eingabewert 15 statt -2 ---
$ws.Range("E10").Value = 15

# --- Row 11: Eingabewert 0 statt Text "f" ---
$ws.Range("E11").Value = 0

# New shared strings must be (re)created in this exact order so the packed
# sharedStrings table matches the canonical ordering after the old, now
# dangling strings are pruned:
#   1) row 11 Ergebnis-Text ("...0 nicht moeglich...")
$ws.Range("F11").Value = '"Übersetzung der Zahl 0 nicht möglich. Version: 1.0"'

#   2) row 18 Aequivalenzklasse Text (gekuerzt)
$ws.Range("D18").Value = "ÄK_pos_1; 0 < x < 11"

#   3) row 19 Aequivalenzklasse Text (neu formatiert) + Eingabewert -3
$ws.Range("D19").Value = "ÄK_neg_1; x < 0"
$ws.Range("E19").Value = -3

#   4) row 20 (vorher leer) - neue Aequivalenzklasse, Eingabewert 15, Ergebnis "neg"
$ws.Range("D20").Value = "ÄK_neg_2; x > 10"
$ws.Range("E20").Value = 15
$ws.Range("F20").Value = "neg"

#   5) row 9 Ergebnis-Text ("...-3 nicht moeglich...")
$ws.Range("F9").Value = '"Übersetzung der Zahl -3 nicht möglich. Version: 1.0"'

#   6) row 10 Ergebnis-Text ("...15 nicht moeglich...")
$ws.Range("F10").Value = '"Übersetzung der Zahl 15 nicht möglich. Version: 1.0"'

# Aktive Selektion auf F12 setzen (UML-Diagramm-Anpassung / letzte Zelle)
$ws.Range("F12").Select()
